$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title text: October -> November ---
$ws.Range("A2").Value = "by State, by Sector, November 2016 and 2015 (Thousand Megawatthours)"

# --- Column header dates (row 6): October 2016/2015 -> November 2016/2015 ---
# Force text format first so Excel does not auto-convert "November 2016" into a date serial.
$headerCells = @("B6","C6","E6","F6","G6","H6","I6","J6","K6","L6")
foreach ($cellref in $headerCells) {
    $ws.Range($cellref).NumberFormat = "@"
}
$ws.Range("B6").Value = "November 2016"
$ws.Range("C6").Value = "November 2015"
$ws.Range("E6").Value = "November 2016"
$ws.Range("F6").Value = "November 2015"
$ws.Range("G6").Value = "November 2016"
$ws.Range("H6").Value = "November 2015"
$ws.Range("I6").Value = "November 2016"
$ws.Range("J6").Value = "November 2015"
$ws.Range("K6").Value = "November 2016"
$ws.Range("L6").Value = "November 2015"
foreach ($cellref in $headerCells) {
    $ws.Range($cellref).NumberFormat = "#,##0"
}

# --- Updated data values (November 2016 vs November 2015 release) ---
$ws.Range("C14").Value = 19
$ws.Range("L14").Value = 19
$ws.Range("B15").Value = "NM"
$ws.Range("C15").Value = 6
$ws.Range("D15").Value = "NM"
$ws.Range("K15").Value = "NM"
$ws.Range("L15").Value = 6
$ws.Range("C17").Value = 13
$ws.Range("L17").Value = 13
$ws.Range("B18").Value = 135
$ws.Range("C18").Value = 317
$ws.Range("D18").Value = -0.575
$ws.Range("E18").Value = 45
$ws.Range("F18").Value = 196
$ws.Range("G18").Value = 76
$ws.Range("H18").Value = 100
$ws.Range("L18").Value = 20
$ws.Range("C20").Value = 179
$ws.Range("F20").Value = 179
$ws.Range("B21").Value = 48
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 0.91
$ws.Range("E21").Value = 36
$ws.Range("F21").Value = 9
$ws.Range("H21").Value = 4
$ws.Range("B22").Value = 77
$ws.Range("C22").Value = 97
$ws.Range("D22").Value = -0.21
$ws.Range("G22").Value = 76
$ws.Range("H22").Value = 96
$ws.Range("B23").Value = 10
$ws.Range("C23").Value = 16
$ws.Range("D23").Value = -0.365
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = 8
$ws.Range("L23").Value = 8
$ws.Range("C24").Value = 4
$ws.Range("I24").Value = 0.31
$ws.Range("L24").Value = 3
$ws.Range("C25").Value = 4
$ws.Range("I25").Value = 0.31
$ws.Range("L25").Value = 3
$ws.Range("B32").Value = 42
$ws.Range("C32").Value = 108
$ws.Range("D32").Value = -0.611
$ws.Range("E32").Value = 35
$ws.Range("F32").Value = 94
$ws.Range("K32").Value = 7
$ws.Range("L32").Value = 14
$ws.Range("B35").Value = 35
$ws.Range("C35").Value = 94
$ws.Range("D35").Value = -0.628
$ws.Range("E35").Value = 35
$ws.Range("F35").Value = 94
$ws.Range("B36").Value = 7
$ws.Range("C36").Value = 14
$ws.Range("D36").Value = -0.495
$ws.Range("K36").Value = 7
$ws.Range("L36").Value = 14
$ws.Range("B42").Value = 79
$ws.Range("C42").Value = 78
$ws.Range("D42").Value = 0.012
$ws.Range("E42").Value = 79
$ws.Range("F42").Value = 78
$ws.Range("B44").Value = 79
$ws.Range("C44").Value = 78
$ws.Range("D44").Value = 0.012
$ws.Range("E44").Value = 79
$ws.Range("F44").Value = 78
$ws.Range("B47").Value = 461
$ws.Range("C47").Value = 151
$ws.Range("D47").Value = 2.059
$ws.Range("E47").Value = 437
$ws.Range("F47").Value = 121
$ws.Range("L47").Value = 30
$ws.Range("B49").Value = 452
$ws.Range("C49").Value = 141
$ws.Range("D49").Value = 2.214
$ws.Range("E49").Value = 437
$ws.Range("F49").Value = 121
$ws.Range("L49").Value = 20
$ws.Range("B52").Value = 40
$ws.Range("C52").Value = 40
$ws.Range("D52").Value = 0.018
$ws.Range("G52").Value = 40
$ws.Range("H52").Value = 40
$ws.Range("B56").Value = 40
$ws.Range("C56").Value = 40
$ws.Range("D56").Value = 0.018
$ws.Range("G56").Value = 40
$ws.Range("H56").Value = 40
$ws.Range("B68").Value = 781
$ws.Range("C68").Value = 715
$ws.Range("D68").Value = 0.092
$ws.Range("E68").Value = 596
$ws.Range("F68").Value = 490
$ws.Range("G68").Value = 116
$ws.Range("H68").Value = 140
$ws.Range("I68").Value = 0.31
$ws.Range("K68").Value = 69
$ws.Range("L68").Value = 85
